$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new columns at E:F for the gridsearch output
#     "Min_of_Loss_Function" and "convergence_code". This shifts the
#     existing "number_options"/"future" columns from E/F to G/H. ---
$ws.Columns("E:F").Insert()

# --- New header cells. Written via a temporary "=\"text\"" formula and
#     then flattened to a literal value with Copy + PasteSpecial(values),
#     so they land as plain shared-string cells (no leftover formula,
#     no style change) just like the rest of the header row. ---
$ws.Cells.Item(1, 5).Formula = '="Min_of_Loss_Function"'
$ws.Cells.Item(1, 6).Formula = '="convergence_code"'
$hdr = $ws.Range("E1:F1")
$hdr.Copy()
$hdr.PasteSpecial(-4163)
$excel.CutCopyMode = 0

# --- Refreshed gridsearch + parameter results. B:D keep the same
#     semantic columns (kappa/theta/sigma) but the new run's values are
#     stored as text, matching the new E (Min_of_Loss_Function) and
#     F (convergence_code) columns added alongside them. ---
$rowValues = @{
    2 = @("42.8763013530862", "23.6098720539757", "18.775679522497", "3.69523976311823")
    3 = @("45.3357355953877", "29.3196151477721", "23.1406231715724", "7.19183864895222")
    4 = @("65.2616518161135", "33.0501361318493", "34.5027277643546", "7.49889322333684")
    5 = @("55.4283633869522", "33.7547600930658", "29.2313237429204", "7.46244060702908")
    6 = @("56.6554010930514", "34.6599768985326", "28.1071616282124", "6.9518800745399")
    7 = @("66.7167033510985", "43.6430444924148", "37.226956067394",  "9.48874584633205")
    8 = @("59.7786538985285", "44.886193000946",  "38.5993295660477", "16.3485341677696")
}

foreach ($r in 2..8) {
    $vals = $rowValues[$r] + @("0")
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $cell = $ws.Cells.Item($r, 2 + $i)   # columns B..F
        $cell.Formula = '="' + $vals[$i] + '"'
    }
    $rowRng = $ws.Range("B$r`:F$r")
    $rowRng.Copy()
    $rowRng.PasteSpecial(-4163)
    $excel.CutCopyMode = 0
}
